{"js": "// Tata Safari blog-post SEO description update:\n//   \"I will highlight about Tata Safari\"  -> \"I will highlight the Tata Safari\"\n//   \"Read on to find it more.\"            -> \"Read on to find out more.\"\n//\n// Both phrases live in the same paragraph (the \"Description:\" paragraph of\n// the Tata Safari blog post). We locate that paragraph first so the\n// replacements cannot accidentally touch the unrelated \"about\"/\"it\"\n// elsewhere in the document, then do a scoped search + replace for each\n// phrase.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Tata Safari 2023 facelift car, prices\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  // \"about\" -> \"the\"\n  const aboutSearch = target.search(\"I will highlight about Tata Safari\", { matchCase: true });\n  aboutSearch.load(\"text\");\n  await context.sync();\n\n  if (aboutSearch.items.length > 0) {\n    aboutSearch.items[0].insertText(\"I will highlight the Tata Safari\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  // \"it\" -> \"out\"\n  const itSearch = target.search(\"Read on to find it more.\", { matchCase: true });\n  itSearch.load(\"text\");\n  await context.sync();\n\n  if (itSearch.items.length > 0) {\n    itSearch.items[0].insertText(\"Read on to find out more.\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Tata Safari blog-post SEO description update:\n#   \"I will highlight about Tata Safari\"  -> \"I will highlight the Tata Safari\"\n#   \"Read on to find it more.\"            -> \"Read on to find out more.\"\n#\n# Both phrases live in the same paragraph (the \"Description:\" paragraph of\n# the Tata Safari blog post). We locate that paragraph first so the\n# replacements cannot accidentally touch the unrelated \"about\"/\"it\"\n# elsewhere in the document, then run a scoped Find/Replace (wdReplaceOne)\n# for each phrase inside that paragraph's range only.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Tata Safari 2023 facelift car, prices*\") {\n\n        # \"about\" -> \"the\"\n        $rng1 = $p.Range\n        $find1 = $rng1.Find\n        $find1.Text = \"I will highlight about Tata Safari\"\n        $find1.Replacement.Text = \"I will highlight the Tata Safari\"\n        $find1.Execute([ref]$find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n        # \"it\" -> \"out\"\n        $rng2 = $p.Range\n        $find2 = $rng2.Find\n        $find2.Text = \"Read on to find it more.\"\n        $find2.Replacement.Text = \"Read on to find out more.\"\n        $find2.Execute([ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\n        break\n    }\n}\n"}
